$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column for rows 2-6 from 2023-09-01 (45170)
# to 2023-09-05 (45174), matching the serial date values stored in the workbook.
$ws.Range("C2:C6").Value = 45174
